$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "[name=`"Greatmouth Mob`"]  Right, I’m curious. I mean, your title’s 'Plastic', and it’s one of your main selling points...`n"
$ws.Range("C12").Value = "[name=`"??? `"]  Certainly, we have something to look forward to... However, before that, please allow me to have a chat with Mr. 'Plastic' Szewczyk, one-on-one.`n"
$ws.Range("C15").Value = "[name=`"Plastic Knight`"]  ...You 'work on' them, huh? What do you want?`n"
$ws.Range("C17").Value = "[name=`"??? `"]  The Roar Guards Company has requested that you wear its latest Jack 2 model body armor in your one-on-one match, or 'duel' as it is more commonly called, with Maria Nearl, and to prolong the match for as long as possible.`n"
$ws.Range("C38").Value = "[name=`"Bald Marcin`"]  Here, your 'Thorn Tear'. It’s not often you order something this strong. Are you worried about her?`n"
$ws.Range("C50").Value = "[name=`"Greatmouth Mob`"]  You’re in for a treat today at this arena, fully sponsored by the Roar Guards Company! 'Roar Guards, even the winds will bow to you!' `n"
$ws.Range("C81").Value = "[name=`"Greatmouth Mob`"]  His might and economical title are unmatched on today’s roster! Please join me and welcome the 'Plastic', Szewczyk!`n"
$ws.Range("C84").Value = "[name=`"Greatmouth Mob`"]  The 'Plastic' Szewczyk, everyone, and check out his armor from the Roar Guards Company, made with the newest materials! Will the latest in Columbian armor technology prove to be the decisive element in this match?! `n"
$ws.Range("C86").Value = "[name=`"Greatmouth Mob`"]  Don’t forget, 'Roar Guards, even the winds will bow to you!' `n"
$ws.Range("C89").Value = "[name=`"Plastic Knight`"]  'Whislash' Zofia. I thought you’d never set foot in this place again after your injury, lest you humiliate yourself.`n"
$ws.Range("C98").Value = "[name=`"Greatmouth Mob`"]  Oh? Ohh? It looked like Szewczyk talked to someone in the audience briefly before he stepped into the ring. And who would’ve thought? It’s 'Whislash', Zofia! What a surprise we have in our audience today!     `n"
$ws.Range("C111").Value = "[name=`"Maria`"]  My name is Maria Nearl! My family motto is, 'Fear neither hardship nor darkness!'`n"
